$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new row at row 124 (pushes old rows 124..252 down to 125..253)
$ws.Rows.Item(124).Insert()

$ws.Cells.Item(124, 1).Value = 8
$ws.Cells.Item(124, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(124, 3).Value = "Coquimbo"
$ws.Cells.Item(124, 4).Value = "12/24/2021"
$ws.Cells.Item(124, 5).Value = 4
$ws.Cells.Item(124, 6).Value = 100114013
$ws.Cells.Item(124, 7).Value = "Zanahoria"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 800
$ws.Cells.Item(124, 11).Value = 6000
$ws.Cells.Item(124, 12).Value = 6500
$ws.Cells.Item(124, 13).Value = 6250
$ws.Cells.Item(124, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(124, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(124, 16).Value = 312
$ws.Cells.Item(124, 17).Value = 20
$ws.Cells.Item(124, 18).Value = "Hortaliza"

# Insert second new row at row 178 (after the first insert, old row 177 now sits at 178;
# this insert pushes it down to 179 and places the new data at 178)
$ws.Rows.Item(178).Insert()

$ws.Cells.Item(178, 1).Value = 8
$ws.Cells.Item(178, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(178, 3).Value = "Coquimbo"
$ws.Cells.Item(178, 4).Value = "12/23/2021"
$ws.Cells.Item(178, 5).Value = 4
$ws.Cells.Item(178, 6).Value = 100114013
$ws.Cells.Item(178, 7).Value = "Zanahoria"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 600
$ws.Cells.Item(178, 11).Value = 6000
$ws.Cells.Item(178, 12).Value = 7000
$ws.Cells.Item(178, 13).Value = 6500
$ws.Cells.Item(178, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(178, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(178, 16).Value = 325
$ws.Cells.Item(178, 17).Value = 20
$ws.Cells.Item(178, 18).Value = "Hortaliza"
